$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix formatting of "fullRNASEQ" -> "fullRNASeq" in the purpose column (E2:E27)
$ws.Range("E2:E27").Value = "fullRNASeq"
